# Apply the scheduled market-data refresh to the per-job Leve-profit tables.
# Updates currentAveragePrice / currentAveragePriceNQ / currentAveragePriceHQ /
# LevePriceNQ / LevePriceHQ / LeveProfitNQ / LeveProfitHQ (columns H-N) with the
# latest Universalis market snapshot for the affected Leve rows, one worksheet
# per crafting job (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR).

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 70
$ws.Range("H70").Value = 8200.538
$ws.Range("I70").Value = 6071
$ws.Range("J70").Value = 8587.727999999999
$ws.Range("K70").Value = 18213
$ws.Range("L70").Value = 25763.184
$ws.Range("M70").Value = -17943
$ws.Range("N70").Value = -26303.184
# Row 73
$ws.Range("H73").Value = 8200.538
$ws.Range("I73").Value = 6071
$ws.Range("J73").Value = 8587.727999999999
$ws.Range("K73").Value = 18213
$ws.Range("L73").Value = 25763.184
$ws.Range("M73").Value = -17277
$ws.Range("N73").Value = -27635.184
# Row 129
$ws.Range("H129").Value = 2540.6667
$ws.Range("I129").Value = 1316.6
$ws.Range("K129").Value = 3949.8
$ws.Range("M129").Value = 1050.2
# Row 137
$ws.Range("H137").Value = 38376.53
$ws.Range("I137").Value = 62167.8
$ws.Range("J137").Value = 28463.5
$ws.Range("K137").Value = 186503.4
$ws.Range("L137").Value = 85390.5
$ws.Range("M137").Value = -183953.4
$ws.Range("N137").Value = -90490.5
# Row 138
$ws.Range("H138").Value = 26854
$ws.Range("I138").Value = 3015.842
$ws.Range("J138").Value = 45725.875
$ws.Range("K138").Value = 9047.526
$ws.Range("L138").Value = 137177.625
$ws.Range("M138").Value = -3907.526
$ws.Range("N138").Value = -147457.625

$ws = $wb.Worksheets.Item("ARM")
# Row 61
$ws.Range("H61").Value = 6947.5483
$ws.Range("I61").Value = 3833.6956
$ws.Range("K61").Value = 3833.6956
$ws.Range("M61").Value = -3621.6956
# Row 132
$ws.Range("H132").Value = 2243.05
$ws.Range("I132").Value = 2021.2858
$ws.Range("J132").Value = 3230.9092
$ws.Range("K132").Value = 6063.857400000001
$ws.Range("L132").Value = 9692.7276
$ws.Range("M132").Value = -3533.857400000001
$ws.Range("N132").Value = -14752.7276
# Row 136
$ws.Range("H136").Value = 6947.5483
$ws.Range("I136").Value = 3833.6956
$ws.Range("K136").Value = 11501.0868
$ws.Range("M136").Value = -8951.086800000001

$ws = $wb.Worksheets.Item("BSM")
# Row 64
$ws.Range("H64").Value = 1256.3334
$ws.Range("I64").Value = 761.8
$ws.Range("J64").Value = 1874.5
$ws.Range("K64").Value = 761.8
$ws.Range("L64").Value = 1874.5
$ws.Range("M64").Value = -536.8
$ws.Range("N64").Value = -2324.5
# Row 67
$ws.Range("H67").Value = 1256.3334
$ws.Range("I67").Value = 761.8
$ws.Range("J67").Value = 1874.5
$ws.Range("K67").Value = 761.8
$ws.Range("L67").Value = 1874.5
$ws.Range("M67").Value = 18.20000000000005
$ws.Range("N67").Value = -3434.5

$ws = $wb.Worksheets.Item("CRP")
# Row 62
$ws.Range("H62").Value = 3757.5
$ws.Range("I62").Value = 3989.4
$ws.Range("J62").Value = 2598
$ws.Range("K62").Value = 3989.4
$ws.Range("L62").Value = 2598
$ws.Range("M62").Value = -3365.4
$ws.Range("N62").Value = -3846
# Row 65
$ws.Range("H65").Value = 3757.5
$ws.Range("I65").Value = 3989.4
$ws.Range("J65").Value = 2598
$ws.Range("K65").Value = 19947
$ws.Range("L65").Value = 12990
$ws.Range("M65").Value = -16827
$ws.Range("N65").Value = -19230
# Row 99
$ws.Range("H99").Value = 7192.2
$ws.Range("I99").Value = 6305.857
$ws.Range("J99").Value = 9260.333000000001
$ws.Range("K99").Value = 6305.857
$ws.Range("L99").Value = 9260.333000000001
$ws.Range("M99").Value = -4807.857
$ws.Range("N99").Value = -12256.333
# Row 107
$ws.Range("H107").Value = 1285.5
$ws.Range("I107").Value = 1000
$ws.Range("K107").Value = 1000
$ws.Range("M107").Value = 920
# Row 126
$ws.Range("H126").Value = 7192.2
$ws.Range("I126").Value = 6305.857
$ws.Range("J126").Value = 9260.333000000001
$ws.Range("K126").Value = 18917.571
$ws.Range("L126").Value = 27780.999
$ws.Range("M126").Value = -16447.571
$ws.Range("N126").Value = -32720.999
# Row 134
$ws.Range("H134").Value = 1919.0652
$ws.Range("I134").Value = 1788.1794
$ws.Range("J134").Value = 2648.2856
$ws.Range("K134").Value = 5364.5382
$ws.Range("L134").Value = 7944.8568
$ws.Range("M134").Value = -2829.5382
$ws.Range("N134").Value = -13014.8568
# Row 141
$ws.Range("H141").Value = 0
$ws.Range("I141").Value = 0
$ws.Range("J141").Value = 0
$ws.Range("K141").Value = 0
$ws.Range("L141").Value = 0
$ws.Range("M141").ClearContents()
$ws.Range("N141").ClearContents()

$ws = $wb.Worksheets.Item("CUL")
# Row 32
$ws.Range("H32").Value = 3999.7778
$ws.Range("J32").Value = 4399.8
$ws.Range("L32").Value = 13199.4
$ws.Range("N32").Value = -13765.4
# Row 68
$ws.Range("H68").Value = 3410.3914
$ws.Range("J68").Value = 3604.6099
$ws.Range("L68").Value = 10813.8297
$ws.Range("N68").Value = -12435.8297
# Row 71
$ws.Range("H71").Value = 3410.3914
$ws.Range("J71").Value = 3604.6099
$ws.Range("L71").Value = 32441.4891
$ws.Range("N71").Value = -40553.4891
# Row 103
$ws.Range("H103").Value = 920
$ws.Range("I103").Value = 1138
$ws.Range("J103").Value = 48
$ws.Range("K103").Value = 3414
$ws.Range("L103").Value = 144
$ws.Range("M103").Value = -2535
$ws.Range("N103").Value = -1902
# Row 113
$ws.Range("H113").Value = 849.875
$ws.Range("I113").Value = 500
$ws.Range("J113").Value = 899.8570999999999
$ws.Range("K113").Value = 1500
$ws.Range("L113").Value = 2699.5713
$ws.Range("M113").Value = 670
$ws.Range("N113").Value = -7039.5713
# Row 129
$ws.Range("H129").Value = 2672.6667
$ws.Range("I129").Value = 1835.8
$ws.Range("K129").Value = 5507.4
$ws.Range("M129").Value = -507.3999999999996

$ws = $wb.Worksheets.Item("GSM")
# Row 57
$ws.Range("H57").Value = 49997.5
$ws.Range("J57").Value = 49997.5
$ws.Range("L57").Value = 49997.5
$ws.Range("N57").Value = -51637.5
# Row 80
$ws.Range("H80").Value = 3951.2
$ws.Range("J80").Value = 4535.846
$ws.Range("L80").Value = 4535.846
$ws.Range("N80").Value = -6531.846
# Row 83
$ws.Range("H83").Value = 3951.2
$ws.Range("J83").Value = 4535.846
$ws.Range("L83").Value = 22679.23
$ws.Range("N83").Value = -32663.23
# Row 126
$ws.Range("H126").Value = 2750
$ws.Range("J126").Value = 2750
$ws.Range("L126").Value = 8250
$ws.Range("N126").Value = -13190
# Row 132
$ws.Range("H132").Value = 1940.931
$ws.Range("I132").Value = 1831.6786
$ws.Range("J132").Value = 5000
$ws.Range("K132").Value = 5495.0358
$ws.Range("L132").Value = 15000
$ws.Range("M132").Value = -2965.0358
$ws.Range("N132").Value = -20060

$ws = $wb.Worksheets.Item("LTW")
# Row 40
$ws.Range("H40").Value = 3597.5557
$ws.Range("I40").Value = 3063.8333
$ws.Range("K40").Value = 3063.8333
$ws.Range("M40").Value = -2927.8333
# Row 61
$ws.Range("H61").Value = 3916.1667
$ws.Range("I61").Value = 4299.4
$ws.Range("J61").Value = 2000
$ws.Range("K61").Value = 4299.4
$ws.Range("L61").Value = 2000
$ws.Range("M61").Value = -4097.4
$ws.Range("N61").Value = -2404
# Row 82
$ws.Range("H82").Value = 3901.2727
$ws.Range("I82").Value = 3627.4614
$ws.Range("K82").Value = 3627.4614
$ws.Range("M82").Value = -3266.4614
# Row 85
$ws.Range("H85").Value = 3901.2727
$ws.Range("I85").Value = 3627.4614
$ws.Range("K85").Value = 3627.4614
$ws.Range("M85").Value = -2379.4614
# Row 110
$ws.Range("H110").Value = 62618
$ws.Range("J110").Value = 62618
$ws.Range("L110").Value = 62618
$ws.Range("N110").Value = -70798
# Row 113
$ws.Range("H113").Value = 3916.1667
$ws.Range("I113").Value = 4299.4
$ws.Range("J113").Value = 2000
$ws.Range("K113").Value = 4299.4
$ws.Range("L113").Value = 2000
$ws.Range("M113").Value = -2129.4
$ws.Range("N113").Value = -6340
# Row 132
$ws.Range("H132").Value = 5445.8047
$ws.Range("I132").Value = 5112.6875
$ws.Range("J132").Value = 6630.222
$ws.Range("K132").Value = 15338.0625
$ws.Range("L132").Value = 19890.666
$ws.Range("M132").Value = -12808.0625
$ws.Range("N132").Value = -24950.666

$ws = $wb.Worksheets.Item("WVR")
# Row 51
$ws.Range("H51").Value = 85000
$ws.Range("I51").Value = 30000
$ws.Range("K51").Value = 30000
$ws.Range("M51").Value = -29490
# Row 59
$ws.Range("H59").Value = 100000
$ws.Range("J59").Value = 100000
$ws.Range("L59").Value = 100000
$ws.Range("N59").Value = -101476
# Row 122
$ws.Range("H122").Value = 23893.234
$ws.Range("I122").Value = 27110.127
$ws.Range("K122").Value = 81330.38099999999
$ws.Range("M122").Value = -78880.38099999999
# Row 132
$ws.Range("H132").Value = 28417.215
$ws.Range("I132").Value = 49216.4
$ws.Range("J132").Value = 4418.154
$ws.Range("K132").Value = 147649.2
$ws.Range("L132").Value = 13254.462
$ws.Range("M132").Value = -145119.2
$ws.Range("N132").Value = -18314.462
# Row 135
$ws.Range("H135").Value = 138962.5
$ws.Range("J135").Value = 138962.5
$ws.Range("L135").Value = 138962.5
$ws.Range("N135").Value = -149102.5
# Row 136
$ws.Range("H136").Value = 14847.178
$ws.Range("I136").Value = 15062.1455
$ws.Range("K136").Value = 45186.4365
$ws.Range("M136").Value = -42636.4365
